$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 319-320; this shifts the existing rows
# (319..418) down to (321..420), which is exactly the re-basing the
# diff shows (every existing record moves down by 2 rows).
$ws.Rows("319:320").Insert()

# Fill in the two brand-new records (week of 44722) in the freshly
# inserted rows, matching the fixed/constant columns used throughout
# this sheet.
$ws.Range("A319").Value = 8
$ws.Range("B319").Value = "Terminal La Palmera de La Serena"
$ws.Range("C319").Value = "Coquimbo"
$ws.Range("D319").Value = 44722
$ws.Range("E319").Value = 4
$ws.Range("F319").Value = 100112017
$ws.Range("G319").Value = "Apio"
$ws.Range("H319").Value = "Americana (o)"
$ws.Range("I319").Value = "Primera"
$ws.Range("J319").Value = 2500
$ws.Range("K319").Value = 7500
$ws.Range("L319").Value = 8000
$ws.Range("M319").Value = 7750
$ws.Range("N319").Value = "`$/docena de matas"
$ws.Range("O319").Value = "Provincia del Elquí"
$ws.Range("P319").Value = 1292
$ws.Range("Q319").Value = 6
$ws.Range("R319").Value = "Hortaliza"

$ws.Range("A320").Value = 8
$ws.Range("B320").Value = "Terminal La Palmera de La Serena"
$ws.Range("C320").Value = "Coquimbo"
$ws.Range("D320").Value = 44722
$ws.Range("E320").Value = 4
$ws.Range("F320").Value = 100112017
$ws.Range("G320").Value = "Apio"
$ws.Range("H320").Value = "Americana (o)"
$ws.Range("I320").Value = "Segunda"
$ws.Range("J320").Value = 1500
$ws.Range("K320").Value = 6500
$ws.Range("L320").Value = 7000
$ws.Range("M320").Value = 6750
$ws.Range("N320").Value = "`$/docena de matas"
$ws.Range("O320").Value = "Provincia del Elquí"
$ws.Range("P320").Value = 1125
$ws.Range("Q320").Value = 6
$ws.Range("R320").Value = "Hortaliza"
